$d = $word.ActiveDocument

# The document ends with a trailing empty paragraph (w:pPr numId=1 ListParagraph)
# that only holds the hidden "_GoBack" bookmark. The edit:
#   1. Inserts a new plain paragraph "Socks in the Dark" just before it.
#   2. Fills that trailing paragraph with the new problem statement text
#      (keeping the bookmark at the end) and moves it onto its own,
#      independent numbered list (numId 3 instead of numId 1).

$target = $d.Paragraphs.Last

# 1) New heading-style paragraph "Socks in the Dark" ------------------------
$target.Range.InsertParagraphBefore()
$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$headingPara.Range.ListFormat.RemoveNumbers()
$headingPara.Range.Style = "Normal"
$headingPara.Range.Text = "Socks in the Dark"

# 2) Add the problem statement text to the (still) trailing paragraph -------
$target = $d.Paragraphs.Last
$insertionPoint = $target.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertBefore("The Problem: There are 20 socks in a drawer 5 pairs are black, 3 pairs are brown, and 2 pairs are white can only select the socks in the dark and can check them only after a selection has been made.")

# 3) Put that paragraph on its own new numbered list (numId 3) --------------
$newListTemplate = $d.ListTemplates.Add($true)
$target.Range.ListFormat.ApplyListTemplateWithLevel($newListTemplate, $false, 1, $false, 1)
